$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.647.69'
$ws.Range("E2").Value = '  -0.05%  '

$ws.Range("D3").Value = '1.642.78'
$ws.Range("E3").Value = '  +0.72%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").Value = '''214.95'
$ws.Range("E5").Value = '  +0.76%  '

$ws.Range("D6").Value = '''0.505'
$ws.Range("E6").Value = '  +0.99%  '

$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = '''0.0627'
$ws.Range("E9").Value = '  +0.78%  '

$ws.Range("D10").Value = '''19.25'
$ws.Range("E10").Value = '  +0.18%  '

$ws.Range("D11").Value = '''0.0842'
$ws.Range("E11").Value = '  +0.17%  '

$ws.Range("D12").Value = '1.871.63'
$ws.Range("E12").Value = '  +0.71%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.20'
$ws.Range("E13").Value = '  +2.74%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.635.23'
$ws.Range("E14").Value = '  +1.88%  '

$ws.Range("E15").Value = '  +1.12%  '

$ws.Range("E16").Value = '  +3.15%  '

$ws.Range("D17").Value = '26.692.04'

$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("D19").Value = '''216.11'
$ws.Range("E19").Value = '  -1.04%  '

$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("E21").Value = '  +1.69%  '

$ws.Range("E22").Value = '  +2.16%  '

$ws.Range("D23").Value = '''9.56'
$ws.Range("E23").Value = '  +2.06%  '

$ws.Range("D24").Value = '''2.18'
$ws.Range("E24").Value = '  +11.35%  '

$ws.Range("D25").Value = '''145.80'
$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("E26").Value = '  +0.21%  '

$ws.Range("E27").Value = '  -0.36%  '

$ws.Range("D28").Value = '''7.16'
$ws.Range("E28").Value = '  +4.61%  '

$ws.Range("D29").Value = '''15.77'
$ws.Range("E29").Value = '  +1.64%  '

$ws.Range("E30").Value = '  +2.43%  '

$ws.Range("E31").Value = '  +0.10%  '

$ws.Range("D32").Value = '''3.39'
$ws.Range("E32").Value = '  +2.61%  '

$ws.Range("E33").Value = '  +2.06%  '

$ws.Range("D34").Value = '1.274.69'
$ws.Range("E34").Value = '  +4.91%  '

$ws.Range("E35").Value = '  +2.05%  '

$ws.Range("E36").Value = '  +5.72%  '

$ws.Range("E37").Value = '  +0.23%  '

$ws.Range("E38").Value = '  +6.13%  '

$ws.Range("D39").Value = '''0.829'
$ws.Range("E39").Value = '  +3.18%  '

$ws.Range("E40").Value = '  +0.28%  '

$ws.Range("D41").Value = '''0.813'
$ws.Range("E41").Value = '  +2.39%  '

$ws.Range("E42").Value = '  -1.66%  '

$ws.Range("E43").Value = '  +2.01%  '

$ws.Range("D44").Value = '1.782.03'
$ws.Range("E44").Value = '  +0.84%  '

$ws.Range("D45").Value = '''92.82'
$ws.Range("E45").Value = '  +0.12%  '

$ws.Range("D46").Value = '''59.43'
$ws.Range("E46").Value = '  +8.02%  '

$ws.Range("E47").Value = '  +2.55%  '

$ws.Range("E48").Value = '  +0.85%  '

$ws.Range("D49").Value = '''7.79'
$ws.Range("E49").Value = '  +2.79%  '

$ws.Range("D50").Value = '''0.0970'
$ws.Range("E50").Value = '  +3.16%  '

$ws.Range("E51").Value = '  -0.67%  '
